# Re-pull data / push all data / mean calculation
# Update column F (dSF) values for the rows whose underlying data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = 7
    "F15" = -3
    "F17" = 3
    "F18" = 0
    "F35" = -3
    "F36" = -2
    "F40" = -5
    "F45" = -3
    "F48" = 2
    "F50" = 1
    "F54" = -3
    "F60" = -11
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
